$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Test_1.2")
$ws2.Activate()
$ws2.Range("D22:D25").Select()
